# Natmi following Dr Hou advice:
# Re-run of the LR-pair (Col8a1-Itga2) NATMI export now includes a third
# sending/target cluster ("ECs") in addition to the original "FAPs"/"sCs",
# so the 2x2 cluster-pair table (rows 2-5) becomes a full 3x3 cluster-pair
# table (rows 2-10) with refreshed specificity/expression statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster | B Ligand symbol | C Receptor symbol |
#          D Target cluster  | E..T numeric expression/specificity stats
$data = @(
    @("ECs", "Col8a1", "Itga2", "ECs", 2, 0.6666666666666666, 3.346688666666667, 10.040066, 0.1019109290125008, 0.1019109290125008, 2, 0.6666666666666666, 1.888921, 5.666763, 0.4551532417350329, 0.4551532417350328, 6.321630502928666, 56.894674526358, 0.04638508970826853, 0.04638508970826852),
    @("ECs", "Col8a1", "Itga2", "FAPs", 2, 0.6666666666666666, 3.346688666666667, 10.040066, 0.1019109290125008, 0.1019109290125008, 3, 1, 1.087098333333333, 3.261295, 0.2619465454094788, 0.2619465454094787, 3.638179671718889, 32.74361704547, 0.0266952157942952, 0.02669521579429519),
    @("ECs", "Col8a1", "Itga2", "sCs", 2, 0.6666666666666666, 3.346688666666667, 10.040066, 0.1019109290125008, 0.1019109290125008, 3, 1, 1.174057666666667, 3.522173, 0.2829002128554884, 0.2829002128554884, 3.929205487046445, 35.36284938341799, 0.02883062350993703, 0.02883062350993703),
    @("FAPs", "Col8a1", "Itga2", "ECs", 3, 1, 26.21646566666666, 78.64939699999999, 0.7983247435368442, 0.7983247435368442, 2, 0.6666666666666666, 1.888921, 5.666763, 0.4551532417350329, 0.4551532417350328, 49.52083254354566, 445.6874928919109, 0.3633600949780834, 0.3633600949780833),
    @("FAPs", "Col8a1", "Itga2", "FAPs", 3, 1, 26.21646566666666, 78.64939699999999, 0.7983247435368442, 0.7983247435368442, 3, 1, 1.087098333333333, 3.261295, 0.2619465454094788, 0.2619465454094787, 28.49987613212388, 256.498885189115, 0.2091184086843844, 0.2091184086843844),
    @("FAPs", "Col8a1", "Itga2", "sCs", 3, 1, 26.21646566666666, 78.64939699999999, 0.7983247435368442, 0.7983247435368442, 3, 1, 1.174057666666667, 3.522173, 0.2829002128554884, 0.2829002128554884, 30.77964250885344, 277.016782579681, 0.2258462398743764, 0.2258462398743764),
    @("sCs", "Col8a1", "Itga2", "ECs", 3, 1, 3.276195666666666, 9.828586999999999, 0.09976432745065499, 0.09976432745065499, 2, 0.6666666666666666, 1.888921, 5.666763, 0.4551532417350329, 0.4551532417350328, 6.188474794875665, 55.69627315388099, 0.04540805704868094, 0.04540805704868094),
    @("sCs", "Col8a1", "Itga2", "FAPs", 3, 1, 3.276195666666666, 9.828586999999999, 0.09976432745065499, 0.09976432745065499, 3, 1, 1.087098333333333, 3.261295, 0.2619465454094788, 0.2619465454094787, 3.561546848907222, 32.05392164016499, 0.0261329209307991, 0.0261329209307991),
    @("sCs", "Col8a1", "Itga2", "sCs", 3, 1, 3.276195666666666, 9.828586999999999, 0.09976432745065499, 0.09976432745065499, 3, 1, 1.174057666666667, 3.522173, 0.2829002128554884, 0.2829002128554884, 3.846442639950111, 34.617983759551, 0.02822334947117494, 0.02822334947117494),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $col = 1 + $j
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}

# Sheet now spans the new 9-row data block (A1:T10).
